$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the width of column B so the newly inserted column C can match it
$bWidth = $ws.Columns("B:B").ColumnWidth

# Insert a new column before the current column C (Investment Date *)
$ws.Columns("C:C").Insert()

# Match the new column's width to column B's width
$ws.Columns("C:C").ColumnWidth = $bWidth

# Header for the new column
$ws.Range("C1").Value = "Pan *"

# Data values for the new column
$ws.Range("C2").Value = "BUHNXDFEA7"
$ws.Range("C3").Value = "BUHNXDFEA7"

# Match the final selection/active cell seen in the saved workbook
$ws.Range("C4").Select() | Out-Null
